$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.384.48'
$ws.Range("E2").Value = '  +2.03%  '
$ws.Range("D3").Value = '2.347.66'
$ws.Range("E3").Value = '  +0.72%  '
$ws.Range("E4").Value = '  -0.20%  '
$ws.Range("D5").Value = "'543.41"
$ws.Range("E5").Value = '  +1.90%  '
$ws.Range("D6").Value = "'135.31"
$ws.Range("E6").Value = '  +1.74%  '
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = '  +0.63%  '
$ws.Range("D8").Value = "'0.563"
$ws.Range("E8").Value = '  +5.19%  '
$ws.Range("E9").Value = '  +0.37%  '
$ws.Range("D10").Value = "'5.65"
$ws.Range("E10").Value = '  +6.20%  '
$ws.Range("D11").Value = "'0.153"
$ws.Range("E11").Value = '  -0.37%  '
$ws.Range("D12").Value = "'0.357"
$ws.Range("E12").Value = '  +3.50%  '
$ws.Range("D13").Value = "'23.85"
$ws.Range("E13").Value = '  +1.13%  '
$ws.Range("D14").Value = '2.764.61'
$ws.Range("E14").Value = '  +0.45%  '
$ws.Range("D15").Value = '58.311.54'
$ws.Range("E15").Value = '  +1.85%  '
$ws.Range("E16").Value = '  +0.24%  '
$ws.Range("D17").Value = '2.337.06'
$ws.Range("E17").Value = '  -0.47%  '
$ws.Range("D18").Value = "'10.74"
$ws.Range("E18").Value = '  +2.64%  '
$ws.Range("D19").Value = "'334.29"
$ws.Range("E19").Value = '  -2.17%  '
$ws.Range("D20").Value = "'4.27"
$ws.Range("E20").Value = '  +2.17%  '
$ws.Range("D21").Value = "'6.67"
$ws.Range("E21").Value = '  -3.48%  '
$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = '  +0.11%  '
$ws.Range("D23").Value = "'5.63"
$ws.Range("E23").Value = '  +1.00%  '
$ws.Range("E24").Value = '  +1.48%  '
$ws.Range("D25").Value = "'0.169"
$ws.Range("E25").Value = '  +1.92%  '
$ws.Range("D26").Value = "'8.53"
$ws.Range("E26").Value = '  -4.00%  '
$ws.Range("D27").Value = "'0.999"
$ws.Range("E27").Value = '  +0.88%  '
$ws.Range("E28").Value = '  +5.62%  '
$ws.Range("E29").Value = '  +1.96%  '
$ws.Range("D30").Value = "'170.72"
$ws.Range("E30").Value = '  +0.51%  '
$ws.Range("D31").Value = '0.0₃0739'
$ws.Range("E31").Value = '  +1.42%  '
$ws.Range("D32").Value = "'6.13"
$ws.Range("E32").Value = '  -0.11%  '
$ws.Range("E33").Value = '  +12.99%  '
$ws.Range("D34").Value = "'18.45"
$ws.Range("E34").Value = '  -0.44%  '
$ws.Range("E35").Value = '  +0.06%  '
$ws.Range("D36").Value = "'4.27"
$ws.Range("E36").Value = '  +6.23%  '
$ws.Range("E37").Value = '  +0.79%  '
$ws.Range("E38").Value = '  -1.73%  '
$ws.Range("D39").Value = "'1.65"
$ws.Range("E39").Value = '  +4.27%  '
$ws.Range("D40").Value = "'39.15"
$ws.Range("E40").Value = '  +0.49%  '
$ws.Range("D41").Value = "'142.69"
$ws.Range("E41").Value = '  -3.45%  '
$ws.Range("D42").Value = "'3.66"
$ws.Range("E42").Value = '  +1.85%  '
$ws.Range("B43").Value = 'Bittensor'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D43").Value = "'289.52"
$ws.Range("E43").Value = '  +1.07%  '
$ws.Range("B44").Value = 'PolygonEcosystemToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D44").Value = "'0.376"
$ws.Range("E44").Value = '  -0.31%  '
$ws.Range("D45").Value = "'0.0940"
$ws.Range("E45").Value = '  +0.91%  '
$ws.Range("D46").Value = "'19.21"
$ws.Range("E46").Value = '  +2.16%  '
$ws.Range("E47").Value = '  -0.11%  '
$ws.Range("E48").Value = '  +0.44%  '
$ws.Range("E49").Value = '  +0.80%  '
$ws.Range("D50").Value = "'0.385"
$ws.Range("E50").Value = '  +1.88%  '
$ws.Range("D51").Value = "'17.49"
$ws.Range("E51").Value = '  +0.23%  '
